$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-20 Thursday" "2024-06-21 Friday"

Replace-Text "498×5=2490" "355×8=2840"
Replace-Text "135×8=1080" "673×7=4711"
Replace-Text "890×7=6230" "573×6=3438"
Replace-Text "418×9=3762" "230×2=460"
Replace-Text "385×9=3465" "638×8=5104"
Replace-Text "373×7=2611" "216×4=864"
Replace-Text "685×4=2740" "599×3=1797"
Replace-Text "238×9=2142" "824×2=1648"
Replace-Text "987×5=4935" "787×8=6296"
Replace-Text "491×5=2455" "815×4=3260"
Replace-Text "294×6=1764" "398×2=796"
Replace-Text "334×4=1336" "271×7=1897"
Replace-Text "297×2=594" "943×9=8487"
Replace-Text "545×5=2725" "762×6=4572"
Replace-Text "109×3=327" "406×7=2842"
Replace-Text "287×7=2009" "494×8=3952"
Replace-Text "210×3=630" "924×9=8316"
Replace-Text "892×9=8028" "307×9=2763"
Replace-Text "598×2=1196" "239×2=478"
Replace-Text "413×9=3717" "846×2=1692"
Replace-Text "895×2=1790" "589×2=1178"
Replace-Text "309×9=2781" "740×9=6660"
Replace-Text "365×6=2190" "460×4=1840"
Replace-Text "544×8=4352" "439×4=1756"
Replace-Text "508×6=3048" "391×2=782"
